$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2 to hold the latest circular entry; this
# shifts all the existing price rows down by one.
$ws.Rows.Item(2).Insert()

# Copy the (now shifted) row 3 formatting into the new row 2 so the
# cell styles match the rest of the data rows.
$ws.Range("A3:F3").Copy()
$ws.Range("A2:F2").PasteSpecial(-4122)

# Populate the new row with the latest circular data.
$ws.Range("A2").Value() = 18
$ws.Range("B2").Value() = "ALUMINIUM INGOT"
$ws.Range("C2").Value() = "IE07"
$ws.Range("D2").Value() = 281.95
$ws.Range("E2").Value() = "19-11-2025"
$ws.Range("F2").Value() = "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-19-11-2025.pdf"

# The row insert does not carry the existing hyperlinks down with it, so
# rebuild the whole hyperlink collection from scratch in row order. This
# naturally reassigns rId1..rId18 to F2..F19 respectively.
$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("F2"), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-19-11-2025.pdf")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-07-11-2025.pdf")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-01-11-2025.pdf")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-30-10-2025.pdf")
$ws.Hyperlinks.Add($ws.Range("F6"), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-25-10-2025.pdf")
$ws.Hyperlinks.Add($ws.Range("F7"), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-17-10-2025.pdf")
$ws.Hyperlinks.Add($ws.Range("F8"), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-14-10-2025.pdf")
$ws.Hyperlinks.Add($ws.Range("F9"), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-09-10-2025.pdf")
$ws.Hyperlinks.Add($ws.Range("F10"), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-01-10-2025.pdf")
$ws.Hyperlinks.Add($ws.Range("F11"), "https://nalcoindia.com/wp-content/uploads/2025/09/INGOT-30-09-2025.pdf")
$ws.Hyperlinks.Add($ws.Range("F12"), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-25-09-2025.pdf")
$ws.Hyperlinks.Add($ws.Range("F13"), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-20-09-2025.pdf")
$ws.Hyperlinks.Add($ws.Range("F14"), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-17-09-2025.pdf")
$ws.Hyperlinks.Add($ws.Range("F15"), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-01-09-2025.pdf")
$ws.Hyperlinks.Add($ws.Range("F16"), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-28-08-2025.pdf")
$ws.Hyperlinks.Add($ws.Range("F17"), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-21-08-2025.pdf")
$ws.Hyperlinks.Add($ws.Range("F18"), "https://nalcoindia.com/wp-content/uploads/2025/08/INGOT-15-08-2025.pdf")
$ws.Hyperlinks.Add($ws.Range("F19"), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-07-08-2025.pdf")

# Adding hyperlinks re-styles the touched cells with a hyperlink font; put
# the plain data-row style back by pasting formatting from a same-style
# neighbour cell in each row.
for ($r = 2; $r -le 19; $r++) {
    $ws.Range("A$r").Copy()
    $ws.Range("F$r").PasteSpecial(-4122)
}
